$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2: update rpc-reply message-id UUID
$f2 = $ws.Range("F2").Value()
$f2 = $f2.Replace("342357b8-bbf8-4a8a-ad58-fd4ab34e0064", "520e4f13-a9d2-4765-af29-ec9aa9f1967e")
$ws.Range("F2").Value = $f2

# H2: update edit-config response message-id, commit response message-id, and flow-id
$h2 = $ws.Range("H2").Value()
$h2 = $h2.Replace("a0600154-baa8-40ad-b48d-1cfb729a6b77", "de63be04-1a6f-44f1-bcff-3adc4ce2c380")
$h2 = $h2.Replace("811238be-b3f5-4f71-8f1c-7ae08adfe7e1", "588f1240-da42-42f6-aff2-c30b880ac0c8")
$h2 = $h2.Replace('nc-ext:flow-id="75"', 'nc-ext:flow-id="239"')
$ws.Range("H2").Value = $h2

# I2: update rpc-reply message-id UUID
$i2 = $ws.Range("I2").Value()
$i2 = $i2.Replace("98824db2-0f57-4825-8385-8c8fd36a6c74", "8b4d38b1-7b84-4915-8a4d-64237a3d48b6")
$ws.Range("I2").Value = $i2
